$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old (buggy) values for the whole used range before re-writing the
# refreshed naive-forecaster output below.
$ws.Range("A1:BA24").ClearContents()

# Remove the two trailing rows that fall out of the new date range (2029-12-31, 2030-12-31)
$ws.Rows("23:24").Delete()

# Remove the trailing column (oldest forecast vintage ages out)
$ws.Columns("BA").Delete()

# Re-populate every surviving cell with its refreshed forecaster value
# (the naive-forecaster bugfix changed the computed values, not just a shift)
$data = @"
1	2	39583
1	3	39765
1	4	39948
1	5	40130
1	6	40310
1	7	40494
1	8	40676
1	9	40862
1	10	41044
1	11	41228
1	12	41409
1	13	41592
1	14	41774
1	15	41957
1	16	42137
1	17	42321
1	18	42503
1	19	42689
1	20	42867
1	21	43053
1	22	43145
1	23	43235
1	24	43326
1	25	43418
1	26	43510
1	27	43600
1	28	43691
1	29	43783
1	30	43875
1	31	43966
1	32	44068
1	33	44159
1	34	44251
1	35	44341
1	36	44432
1	37	44525
1	38	44617
1	39	44706
1	40	44798
1	41	44890
1	42	44981
1	43	45071
1	44	45163
1	45	45254
1	46	45345
1	47	45436
1	48	45534
1	49	45618
1	50	45713
1	51	45800
1	52	45891
2	1	39813
3	1	40178
3	5	-0.5555135891318952
3	6	-0.5555135891318952
3	7	-0.5555135891318952
3	8	-0.5555135891318952
3	9	-0.5555135891318952
3	10	-0.5555135891318952
3	11	-0.5555135891318952
3	12	-0.5555135891318952
3	13	-0.5555135891318952
3	14	-0.5555135891318952
3	15	-0.5555135891318952
3	16	-0.5555135891318952
3	17	-0.5555135891318952
3	18	-0.5555135891318952
3	19	-0.5555135891318952
3	20	-0.5555135891318952
3	21	-0.5555135891318952
3	22	-0.5555135891318952
3	23	-0.5555135891318952
3	24	-0.5555135891318952
3	25	-0.5555135891318952
3	26	-0.5555135891318952
3	27	-0.5555135891318952
3	28	-0.5555135891318952
3	29	-0.5555135891318952
3	30	-0.5555135891318952
3	31	-0.5555135891318952
3	32	-0.5555135891318952
3	33	-0.5555135891318952
3	34	-0.5555135891318952
3	35	-0.5555135891318952
3	36	-0.5555135891318952
3	37	-0.5555135891318952
3	38	-0.5555135891318952
3	39	-0.5555135891318952
3	40	-0.5555135891318952
3	41	-0.5555135891318952
3	42	-0.5555135891318952
3	43	-0.5555135891318952
3	44	-0.5555135891318952
3	45	-0.5555135891318952
3	46	-0.5555135891318952
3	47	-0.5555135891318952
3	48	-0.5555135891318952
3	49	-0.5555135891318952
3	50	-0.5555135891318952
3	51	-0.5555135891318952
3	52	-0.5555135891318952
4	1	40543
4	7	0.8442071301477228
4	8	0.8442071301477228
4	9	0.8442071301477228
4	10	0.8442071301477228
4	11	0.8442071301477228
4	12	0.8442071301477228
4	13	0.8442071301477228
4	14	0.8442071301477228
4	15	0.8442071301477228
4	16	0.8442071301477228
4	17	0.8442071301477228
4	18	0.8442071301477228
4	19	0.8442071301477228
4	20	0.8442071301477228
4	21	0.8442071301477228
4	22	0.8442071301477228
4	23	0.8442071301477228
4	24	0.8442071301477228
4	25	0.8442071301477228
4	26	0.8442071301477228
4	27	0.8442071301477228
4	28	0.8442071301477228
4	29	0.8442071301477228
4	30	0.8442071301477228
4	31	0.8442071301477228
4	32	0.8442071301477228
4	33	0.8442071301477228
4	34	0.8442071301477228
4	35	0.8442071301477228
4	36	0.8442071301477228
4	37	0.8442071301477228
4	38	0.8442071301477228
4	39	0.8442071301477228
4	40	0.8442071301477228
4	41	0.8442071301477228
4	42	0.8442071301477228
4	43	0.8442071301477228
4	44	0.8442071301477228
4	45	0.8442071301477228
4	46	0.8442071301477228
4	47	0.8442071301477228
4	48	0.8442071301477228
4	49	0.8442071301477228
4	50	0.8442071301477228
4	51	0.8442071301477228
4	52	0.8442071301477228
5	1	40908
5	9	1.122475521884692
5	10	1.122475521884692
5	11	1.122475521884692
5	12	1.122475521884692
5	13	1.122475521884692
5	14	1.122475521884692
5	15	1.122475521884692
5	16	1.122475521884692
5	17	1.122475521884692
5	18	1.122475521884692
5	19	1.122475521884692
5	20	1.122475521884692
5	21	1.122475521884692
5	22	1.122475521884692
5	23	1.122475521884692
5	24	1.122475521884692
5	25	1.122475521884692
5	26	1.122475521884692
5	27	1.122475521884692
5	28	1.122475521884692
5	29	1.122475521884692
5	30	1.122475521884692
5	31	1.122475521884692
5	32	1.122475521884692
5	33	1.122475521884692
5	34	1.122475521884692
5	35	1.122475521884692
5	36	1.122475521884692
5	37	1.122475521884692
5	38	1.122475521884692
5	39	1.122475521884692
5	40	1.122475521884692
5	41	1.122475521884692
5	42	1.122475521884692
5	43	1.122475521884692
5	44	1.122475521884692
5	45	1.122475521884692
5	46	1.122475521884692
5	47	1.122475521884692
5	48	1.122475521884692
5	49	1.122475521884692
5	50	1.122475521884692
5	51	1.122475521884692
5	52	1.122475521884692
6	1	41274
6	11	-0.578174579726376
6	12	-0.578174579726376
6	13	-0.578174579726376
6	14	-0.578174579726376
6	15	-0.578174579726376
6	16	-0.578174579726376
6	17	-0.578174579726376
6	18	-0.578174579726376
6	19	-0.578174579726376
6	20	-0.578174579726376
6	21	-0.578174579726376
6	22	-0.578174579726376
6	23	-0.578174579726376
6	24	-0.578174579726376
6	25	-0.578174579726376
6	26	-0.578174579726376
6	27	-0.578174579726376
6	28	-0.578174579726376
6	29	-0.578174579726376
6	30	-0.578174579726376
6	31	-0.578174579726376
6	32	-0.578174579726376
6	33	-0.578174579726376
6	34	-0.578174579726376
6	35	-0.578174579726376
6	36	-0.578174579726376
6	37	-0.578174579726376
6	38	-0.578174579726376
6	39	-0.578174579726376
6	40	-0.578174579726376
6	41	-0.578174579726376
6	42	-0.578174579726376
6	43	-0.578174579726376
6	44	-0.578174579726376
6	45	-0.578174579726376
6	46	-0.578174579726376
6	47	-0.578174579726376
6	48	-0.578174579726376
6	49	-0.578174579726376
6	50	-0.578174579726376
6	51	-0.578174579726376
6	52	-0.578174579726376
7	1	41639
7	11	-0.3496173419443749
7	12	-0.6616365666142765
7	13	-0.7492845378401558
7	14	-0.7492845378401558
7	15	-0.7492845378401558
7	16	-0.7492845378401558
7	17	-0.7492845378401558
7	18	-0.7492845378401558
7	19	-0.7492845378401558
7	20	-0.7492845378401558
7	21	-0.7492845378401558
7	22	-0.7492845378401558
7	23	-0.7492845378401558
7	24	-0.7492845378401558
7	25	-0.7492845378401558
7	26	-0.7492845378401558
7	27	-0.7492845378401558
7	28	-0.7492845378401558
7	29	-0.7492845378401558
7	30	-0.7492845378401558
7	31	-0.7492845378401558
7	32	-0.7492845378401558
7	33	-0.7492845378401558
7	34	-0.7492845378401558
7	35	-0.7492845378401558
7	36	-0.7492845378401558
7	37	-0.7492845378401558
7	38	-0.7492845378401558
7	39	-0.7492845378401558
7	40	-0.7492845378401558
7	41	-0.7492845378401558
7	42	-0.7492845378401558
7	43	-0.7492845378401558
7	44	-0.7492845378401558
7	45	-0.7492845378401558
7	46	-0.7492845378401558
7	47	-0.7492845378401558
7	48	-0.7492845378401558
7	49	-0.7492845378401558
7	50	-0.7492845378401558
7	51	-0.7492845378401558
7	52	-0.7492845378401558
8	1	42004
8	11	0.1626813993622633
8	12	0.2197847717222867
8	13	0.05500386022236903
8	14	0.2184978785563896
8	15	0.2751437421933511
8	16	0.2751437421933511
8	17	0.2751437421933511
8	18	0.2751437421933511
8	19	0.2751437421933511
8	20	0.2751437421933511
8	21	0.2751437421933511
8	22	0.2751437421933511
8	23	0.2751437421933511
8	24	0.2751437421933511
8	25	0.2751437421933511
8	26	0.2751437421933511
8	27	0.2751437421933511
8	28	0.2751437421933511
8	29	0.2751437421933511
8	30	0.2751437421933511
8	31	0.2751437421933511
8	32	0.2751437421933511
8	33	0.2751437421933511
8	34	0.2751437421933511
8	35	0.2751437421933511
8	36	0.2751437421933511
8	37	0.2751437421933511
8	38	0.2751437421933511
8	39	0.2751437421933511
8	40	0.2751437421933511
8	41	0.2751437421933511
8	42	0.2751437421933511
8	43	0.2751437421933511
8	44	0.2751437421933511
8	45	0.2751437421933511
8	46	0.2751437421933511
8	47	0.2751437421933511
8	48	0.2751437421933511
8	49	0.2751437421933511
8	50	0.2751437421933511
8	51	0.2751437421933511
8	52	0.2751437421933511
9	1	42369
9	12	0.09544618126309246
9	13	0.05909805310246874
9	14	0.04195831742983547
9	15	0.07916875696107883
9	16	-0.01790997771649039
9	17	0.07468705617190707
9	18	0.07468705617190707
9	19	0.07468705617190707
9	20	0.07468705617190707
9	21	0.07468705617190707
9	22	0.07468705617190707
9	23	0.07468705617190707
9	24	0.07468705617190707
9	25	0.07468705617190707
9	26	0.07468705617190707
9	27	0.07468705617190707
9	28	0.07468705617190707
9	29	0.07468705617190707
9	30	0.07468705617190707
9	31	0.07468705617190707
9	32	0.07468705617190707
9	33	0.07468705617190707
9	34	0.07468705617190707
9	35	0.07468705617190707
9	36	0.07468705617190707
9	37	0.07468705617190707
9	38	0.07468705617190707
9	39	0.07468705617190707
9	40	0.07468705617190707
9	41	0.07468705617190707
9	42	0.07468705617190707
9	43	0.07468705617190707
9	44	0.07468705617190707
9	45	0.07468705617190707
9	46	0.07468705617190707
9	47	0.07468705617190707
9	48	0.07468705617190707
9	49	0.07468705617190707
9	50	0.07468705617190707
9	51	0.07468705617190707
9	52	0.07468705617190707
10	1	42735
10	14	0.08086314912676418
10	15	0.09103564879091586
10	16	0.1671491311400208
10	17	0.1656566557188155
10	18	-0.08711135105702317
10	19	-0.05493014849097255
10	20	-0.05493014849097255
10	21	-0.05493014849097255
10	22	-0.05493014849097255
10	23	-0.05493014849097255
10	24	-0.05493014849097255
10	25	-0.05493014849097255
10	26	-0.05493014849097255
10	27	-0.05493014849097255
10	28	-0.05493014849097255
10	29	-0.05493014849097255
10	30	-0.05493014849097255
10	31	-0.05493014849097255
10	32	-0.05493014849097255
10	33	-0.05493014849097255
10	34	-0.05493014849097255
10	35	-0.05493014849097255
10	36	-0.05493014849097255
10	37	-0.05493014849097255
10	38	-0.05493014849097255
10	39	-0.05493014849097255
10	40	-0.05493014849097255
10	41	-0.05493014849097255
10	42	-0.05493014849097255
10	43	-0.05493014849097255
10	44	-0.05493014849097255
10	45	-0.05493014849097255
10	46	-0.05493014849097255
10	47	-0.05493014849097255
10	48	-0.05493014849097255
10	49	-0.05493014849097255
10	50	-0.05493014849097255
10	51	-0.05493014849097255
10	52	-0.05493014849097255
11	1	43100
11	16	0.09001325883963851
11	17	0.09512483792448734
11	18	0.1158714888162216
11	19	0.2044493994367125
11	20	0.1761917659537371
11	21	0.2820931576894115
11	22	0.2820931576894115
11	23	0.2820931576894115
11	24	0.2820931576894115
11	25	0.2820931576894115
11	26	0.2820931576894115
11	27	0.2820931576894115
11	28	0.2820931576894115
11	29	0.2820931576894115
11	30	0.2820931576894115
11	31	0.2820931576894115
11	32	0.2820931576894115
11	33	0.2820931576894115
11	34	0.2820931576894115
11	35	0.2820931576894115
11	36	0.2820931576894115
11	37	0.2820931576894115
11	38	0.2820931576894115
11	39	0.2820931576894115
11	40	0.2820931576894115
11	41	0.2820931576894115
11	42	0.2820931576894115
11	43	0.2820931576894115
11	44	0.2820931576894115
11	45	0.2820931576894115
11	46	0.2820931576894115
11	47	0.2820931576894115
11	48	0.2820931576894115
11	49	0.2820931576894115
11	50	0.2820931576894115
11	51	0.2820931576894115
11	52	0.2820931576894115
12	1	43465
12	18	0.07295307304728826
12	19	0.08337410701473313
12	20	0.1078587431702305
12	21	0.2685680645708288
12	22	0.2803378563356329
12	23	0.3524405906205841
12	24	0.1415392254179304
12	25	0.2343541283920114
12	26	0.2343541283920114
12	27	0.2343541283920114
12	28	0.2343541283920114
12	29	0.2343541283920114
12	30	0.2343541283920114
12	31	0.2343541283920114
12	32	0.2343541283920114
12	33	0.2343541283920114
12	34	0.2343541283920114
12	35	0.2343541283920114
12	36	0.2343541283920114
12	37	0.2343541283920114
12	38	0.2343541283920114
12	39	0.2343541283920114
12	40	0.2343541283920114
12	41	0.2343541283920114
12	42	0.2343541283920114
12	43	0.2343541283920114
12	44	0.2343541283920114
12	45	0.2343541283920114
12	46	0.2343541283920114
12	47	0.2343541283920114
12	48	0.2343541283920114
12	49	0.2343541283920114
12	50	0.2343541283920114
12	51	0.2343541283920114
12	52	0.2343541283920114
13	1	43830
13	20	0.08648845420498041
13	21	0.1068584531153549
13	22	0.1089201880626334
13	23	0.1531204771924033
13	24	-0.1904239862803969
13	25	0.2977174885593792
13	26	0.3424613118119479
13	27	0.07560805834034845
13	28	-0.05069288950212414
13	29	-0.009430310228020211
13	30	-0.009430310228020211
13	31	-0.009430310228020211
13	32	-0.009430310228020211
13	33	-0.009430310228020211
13	34	-0.009430310228020211
13	35	-0.009430310228020211
13	36	-0.009430310228020211
13	37	-0.009430310228020211
13	38	-0.009430310228020211
13	39	-0.009430310228020211
13	40	-0.009430310228020211
13	41	-0.009430310228020211
13	42	-0.009430310228020211
13	43	-0.009430310228020211
13	44	-0.009430310228020211
13	45	-0.009430310228020211
13	46	-0.009430310228020211
13	47	-0.009430310228020211
13	48	-0.009430310228020211
13	49	-0.009430310228020211
13	50	-0.009430310228020211
13	51	-0.009430310228020211
13	52	-0.009430310228020211
14	1	44196
14	23	0.1240365846986169
14	24	0.0862183575680131
14	25	0.1275780827634909
14	26	0.1339497680586277
14	27	-0.01670081902098719
14	28	-0.2043633904885378
14	29	0.0103609600907939
14	30	0.01105513701109562
14	31	-0.5849047489490333
14	32	-2.657403949513992
14	33	-2.657403949513992
14	34	-2.657403949513992
14	35	-2.657403949513992
14	36	-2.657403949513992
14	37	-2.657403949513992
14	38	-2.657403949513992
14	39	-2.657403949513992
14	40	-2.657403949513992
14	41	-2.657403949513992
14	42	-2.657403949513992
14	43	-2.657403949513992
14	44	-2.657403949513992
14	45	-2.657403949513992
14	46	-2.657403949513992
14	47	-2.657403949513992
14	48	-2.657403949513992
14	49	-2.657403949513992
14	50	-2.657403949513992
14	51	-2.657403949513992
14	52	-2.657403949513992
15	1	44561
15	27	0.09988477497939741
15	28	0.08096036838765031
15	29	0.0924908932996793
15	30	0.0925667197466451
15	31	-0.2059746096811033
15	32	-1.40802832891157
15	33	-1.407243743159736
15	34	-0.2500618974080826
15	35	-0.4334047671505248
15	36	-0.3096364143617802
15	37	-0.3096364143617802
15	38	-0.3096364143617802
15	39	-0.3096364143617802
15	40	-0.3096364143617802
15	41	-0.3096364143617802
15	42	-0.3096364143617802
15	43	-0.3096364143617802
15	44	-0.3096364143617802
15	45	-0.3096364143617802
15	46	-0.3096364143617802
15	47	-0.3096364143617802
15	48	-0.3096364143617802
15	49	-0.3096364143617802
15	50	-0.3096364143617802
15	51	-0.3096364143617802
15	52	-0.3096364143617802
16	1	44926
16	31	0.02690490217465147
16	32	-0.3373802571851825
16	33	-0.3318847693681293
16	34	0.07829984441984905
16	35	-0.06335028919957075
16	36	0.2283024244226883
16	37	0.2048390592685578
16	38	-0.0355780787674953
16	39	-0.1663214453978101
16	40	-0.1730430455425092
16	41	-0.1730430455425092
16	42	-0.1730430455425092
16	43	-0.1730430455425092
16	44	-0.1730430455425092
16	45	-0.1730430455425092
16	46	-0.1730430455425092
16	47	-0.1730430455425092
16	48	-0.1730430455425092
16	49	-0.1730430455425092
16	50	-0.1730430455425092
16	51	-0.1730430455425092
16	52	-0.1730430455425092
17	1	45291
17	34	-0.1014535848389841
17	35	-0.1084535449743185
17	36	-0.07323448430569535
17	37	-0.09464543652764057
17	38	-0.1890623092888566
17	39	-0.3613518455741316
17	40	-0.4115424244148125
17	41	0.2152263639657814
17	42	0.6050248749486009
17	43	0.6502606143725664
17	44	0.6376744206510576
17	45	0.6376744206510576
17	46	0.6376744206510576
17	47	0.6376744206510576
17	48	0.6376744206510576
17	49	0.6376744206510576
17	50	0.6376744206510576
17	51	0.6376744206510576
17	52	0.6376744206510576
18	1	45657
18	38	-0.1350767357100935
18	39	-0.2302498790168306
18	40	-0.2817956528829213
18	41	-0.2954722246111707
18	42	-0.1865259660156937
18	43	-0.1691853834640433
18	44	-0.2058599286704377
18	45	-0.06071040501895997
18	46	-0.05036452040672046
18	47	0.1549171986535924
18	48	0.1856341247700399
18	49	0.1856341247700399
18	50	0.1856341247700399
18	51	0.1856341247700399
18	52	0.1856341247700399
19	1	46022
19	42	-0.2229868532896306
19	43	-0.2399087565268632
19	44	-0.2555246202002537
19	45	-0.2428218159789997
19	46	-0.2233431232791294
19	47	-0.03486668218654065
19	48	0.0444032571666142
19	49	0.160714157635633
19	50	0.01130372647704103
19	51	-0.06391119588061711
19	52	-0.09450306168263811
20	1	46387
20	46	-0.214870093455366
20	47	-0.1571058138897907
20	48	-0.1327251081615577
20	49	-0.1040004763365077
20	50	-0.1505299872523014
20	51	-0.2122873162357264
20	52	-0.2824524929558314
21	1	46752
21	50	-0.1156228055083641
21	51	-0.1284476411859137
21	52	-0.1394689437024588
22	1	47118
"@

$lines = $data -split "`n"
foreach ($line in $lines) {
    $t = $line.Trim()
    if ($t -eq "") { continue }
    $parts = $t -split "\s+"
    $r = [int]$parts[0]
    $c = [int]$parts[1]
    $v = [double]$parts[2]
    $ws.Cells.Item($r, $c).Value = $v
}
